# Remove the existing comment ("tag as measurement?") from the document.
# This deletes the w:commentRangeStart / w:commentRangeEnd / w:commentReference
# markers around the word "p" in "...de trois p|asses." and removes the
# corresponding entry from word/comments.xml, while leaving the commented
# text itself untouched.

$d = $word.ActiveDocument

for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments($i).Delete()
}
